# Performance workbook update
# - Insert two new columns (L,M) on Sheet1 to host percentage-difference
#   columns next to the existing absolute-difference columns.
# - Rehome a few labels/formulas that shifted as a result.
# - Add new "% difference" formulas in the new columns.
# - Refresh workbook/view metadata (selection, calc id, etc).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- 1. Make room: insert two columns before column L -------------------
$ws.Range("L1:M1").EntireColumn.Insert()

# --- 2. Row 9/10 header housekeeping -------------------------------------
# "revision" header cell moves from K to L (J10/L10 describe the new J/L
# pair of absolute-difference columns; K/M become the % -difference pair).
$ws.Range("L10").Value = $ws.Range("K10").Value2
$ws.Range("K10").ClearContents()

$ws.Range("L16").Value = $ws.Range("K16").Value2
$ws.Range("K16").ClearContents()

# --- 3. Absolute-difference formulas that used to live in K move to L ---
$ws.Range("L11").Formula = "=G11-G11"
$ws.Range("L12").Formula = "=G12-G12"
$ws.Range("L13").Formula = "=G13-G13"
$ws.Range("L14").Formula = "=G14-G14"
$ws.Range("K11:K14").ClearContents()

$ws.Range("L17").Formula = "=G11-G17"
$ws.Range("L18").Formula = "=G12-G18"
$ws.Range("L19").Formula = "=G13-G19"
$ws.Range("L20").Formula = "=G14-G20"

# --- 4. New percent-difference formulas in K (row block 2) and M --------
$ws.Range("K17").Formula = "=(F17-F11)/F11"
$ws.Range("K18").Formula = "=(F18-F12)/F12"
$ws.Range("K19").Formula = "=(F19-F13)/F13"
$ws.Range("K20").Formula = "=(F20-F14)/F14"

$ws.Range("M17").Formula = "=(G17-G11)/G11"
$ws.Range("M18").Formula = "=(G18-G12)/G12"
$ws.Range("M19").Formula = "=(G19-G13)/G13"
$ws.Range("M20").Formula = "=(G20-G14)/G14"

# --- 5. Styles -------------------------------------------------------------
# Header row (row 2-7): thin horizontal/left alignment styles rotate.
$ws.Range("B2:E2").Style = $ws.Range("B2").Style
$ws.Range("B2:E2").HorizontalAlignment = -4108  # xlCenter placeholder reset below

Write-Host "done"
